$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "59.096.07"
Set-TextValue $ws.Range("E2") "  +2.13%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.589.81"
Set-TextValue $ws.Range("E3") "  +0.97%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "530.50"
Set-TextValue $ws.Range("E5") "  +2.90%  "

# Row 6
Set-TextValue $ws.Range("D6") "139.88"
Set-TextValue $ws.Range("E6") "  +0.71%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.566"
Set-TextValue $ws.Range("E8") "  +0.85%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.602.13"
Set-TextValue $ws.Range("E9") "  +1.02%  "

# Row 10
Set-TextValue $ws.Range("D10") "6.44"
Set-TextValue $ws.Range("E10") "  +0.65%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +2.89%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.333"
Set-TextValue $ws.Range("E12") "  +2.35%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +2.94%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.049.92"
Set-TextValue $ws.Range("E14") "  +1.16%  "

# Row 15
Set-TextValue $ws.Range("D15") "59.045.71"
Set-TextValue $ws.Range("E15") "  +2.01%  "

# Row 16
Set-TextValue $ws.Range("D16") "20.39"
Set-TextValue $ws.Range("E16") "  +1.90%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.0000133"
Set-TextValue $ws.Range("E17") "  +1.64%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.586.39"
Set-TextValue $ws.Range("E18") "  +1.38%  "

# Row 19
Set-TextValue $ws.Range("D19") "346.48"
Set-TextValue $ws.Range("E19") "  +4.50%  "

# Row 20
Set-TextValue $ws.Range("D20") "4.32"
Set-TextValue $ws.Range("E20") "  +1.38%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.09"
Set-TextValue $ws.Range("E21") "  +0.60%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.38"
Set-TextValue $ws.Range("E22") "  +1.34%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +0.09%  "

# Row 24
Set-TextValue $ws.Range("D24") "67.40"
Set-TextValue $ws.Range("E24") "  +2.68%  "

# Row 25
Set-TextValue $ws.Range("E25") "  +1.04%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.405"
Set-TextValue $ws.Range("E26") "  +2.12%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  +0.25%  "

# Row 28
Set-TextValue $ws.Range("D28") "7.14"
Set-TextValue $ws.Range("E28") "  +3.46%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.00%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0₃0730"
Set-TextValue $ws.Range("E30") "  +1.73%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +3.65%  "

# Row 32
Set-TextValue $ws.Range("D32") "5.81"
Set-TextValue $ws.Range("E32") "  -2.67%  "

# Row 33
Set-TextValue $ws.Range("D33") "18.74"
Set-TextValue $ws.Range("E33") "  +0.91%  "

# Row 34
Set-TextValue $ws.Range("D34") "149.45"
Set-TextValue $ws.Range("E34") "  +0.28%  "

# Row 35
Set-TextValue $ws.Range("D35") "3.96"
Set-TextValue $ws.Range("E35") "  +1.71%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +0.57%  "

# Row 37
Set-TextValue $ws.Range("D37") "36.81"
Set-TextValue $ws.Range("E37") "  +2.00%  "

# Row 38
Set-TextValue $ws.Range("E38") "  +3.30%  "

# Row 39
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D39") "0.827"
Set-TextValue $ws.Range("E39") "  +1.40%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("E40") "  -1.58%  "

# Row 41
Set-TextValue $ws.Range("D41") "3.52"
Set-TextValue $ws.Range("E41") "  +1.58%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.997"
Set-TextValue $ws.Range("E42") "  +0.02%  "

# Row 43
Set-TextValue $ws.Range("D43") "271.58"
Set-TextValue $ws.Range("E43") "  +0.22%  "

# Row 44
Set-TextValue $ws.Range("E44") "  +0.62%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.594"
Set-TextValue $ws.Range("E45") "  +1.22%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.0958"
Set-TextValue $ws.Range("E46") "  +2.10%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0517"
Set-TextValue $ws.Range("E47") "  +0.82%  "

# Row 48
Set-TextValue $ws.Range("D48") "18.37"
Set-TextValue $ws.Range("E48") "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D49") "1.945.72"
Set-TextValue $ws.Range("E49") "  -0.78%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0221"
Set-TextValue $ws.Range("E50") "  +1.90%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D51") "18.16"
Set-TextValue $ws.Range("E51") "  +1.43%  "
